$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 4 new data rows right after the current last data row (19)
#    so they become rows 20-23. This pushes the blank spacer rows and
#    the signature/footer block (previously rows 24-25) down to 28-29.
# ------------------------------------------------------------------
$ws.Rows("20:23").Insert() | Out-Null

# Copy the formatting (borders, fills, number formats) of row 19 into
# the 4 newly inserted rows so they look like the rest of the table.
$srcRow = $ws.Range("B19:J19")
foreach ($r in 20..23) {
    $destRow = $ws.Range("B" + $r + ":J" + $r)
    $srcRow.Copy($destRow)
}

# ------------------------------------------------------------------
# 2. Fill in the new period (2509) rows with the same 4 workers that
#    are already listed for period 2508, changing only the period.
# ------------------------------------------------------------------
$ws.Range("C20").Value2 = "1047383546"
$ws.Range("D20").Value2 = "CARLOS MARIO YANES DE LA CRUZ"
$ws.Range("E20").Value2 = "2509"

$ws.Range("C21").Value2 = "1143374181"
$ws.Range("D21").Value2 = "JOSE JESUS ORTEGA CARABALLO"
$ws.Range("E21").Value2 = "2509"

$ws.Range("C22").Value2 = "1049934811"
$ws.Range("D22").Value2 = "JESUS DAVID PATERNINA BARRIOS"
$ws.Range("E22").Value2 = "2509"

$ws.Range("C23").Value2 = "1082887053"
$ws.Range("D23").Value2 = "ELIAS DAVID CHARRIZ VEGA"
$ws.Range("E23").Value2 = "2509"

# ------------------------------------------------------------------
# 3. Center-align the "Periodo Mora" column for every data row so the
#    period values (2508 / 2509) line up consistently.
# ------------------------------------------------------------------
$ws.Range("E16:E23").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 4. Update the summary figures: total overdue value doubles (new
#    period added) and the period count goes from 1 to 2.
# ------------------------------------------------------------------
$ws.Range("E11").Value2 = 508136
$ws.Range("F13").Value2 = 2

Write-Output "Edit complete"
